$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 3-way dipswitches ---
$ws.Range("B7").Value = "3 way dipswitches"
$ws.Range("C7").Value = "adressing for leds"
[void]$ws.Hyperlinks.Add($ws.Range("E7"), "http://nl.farnell.com/omron-electronic-components/a6s3102h/switch-dip-3-way-sealed/dp/1960896")
$ws.Range("D7").Value = 6

# --- Row 8: Ledstrip connector ---
$ws.Range("B8").Value = "Ledstrip connector"
$ws.Range("C8").Value = "to connect strip to board"
[void]$ws.Hyperlinks.Add($ws.Range("E8"), "http://fi.farnell.com/phoenix-contact/ptf-0-3-4-wb-1-8-h/plug-in-connector-pcb-4way-w-wires/dp/2365432?ost=ledstrip&categoryId=700000005017")
$ws.Range("D8").Value = 8

# Hyperlinks.Add re-applies the built-in "Hyperlink" cell style (changing the
# cell's border/alignment formatting). Restore the original look of these
# cells (same formatting already used by the other empty "Shop/Example/Info:"
# cells in the table, e.g. E9) by copying formats only - this keeps the
# underlying cell style index unchanged while leaving the newly-set value and
# hyperlink relationship intact.
$ws.Range("E9").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to E9, matching the saved cursor position.
[void]$ws.Range("E9").Select()
